# Updated symbol list on Tue Jan  3 22:18:05 UTC 2023 with GitHub Actions
#
# Applies updated price/volume figures (and, for rows 10-15, a reshuffled
# coin ordering) to the cryptos worksheet. Columns D (Price) and E
# (Volume(1h)) hold numeric-/percent-looking text that must stay literal
# text (matching the original t="inlineStr" cells), so those assignments
# are prefixed with a leading apostrophe to force Excel to store them as
# text instead of auto-converting to a number/percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.85"
$ws.Range("E2").Value = "'-0.53%"

$ws.Range("D3").Value = "'28.36"
$ws.Range("E3").Value = "'-3.80%"

$ws.Range("D4").Value = "'5.292"
$ws.Range("E4").Value = "'1.84%"

$ws.Range("D5").Value = "'0.05709"
$ws.Range("E5").Value = "'-0.47%"

$ws.Range("E6").Value = "'1.36%"

$ws.Range("D7").Value = "'3.212"
$ws.Range("E7").Value = "'3.61%"

$ws.Range("D8").Value = "'0.8519"
$ws.Range("E8").Value = "'-0.63%"

$ws.Range("D9").Value = "'0.8850"
$ws.Range("E9").Value = "'1.94%"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1391"
$ws.Range("E10").Value = "'1.79%"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07093"
$ws.Range("E11").Value = "'0.26%"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03169"
$ws.Range("E12").Value = "'3.22%"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09227"
$ws.Range("E13").Value = "'-1.69%"

$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001533"
$ws.Range("E14").Value = "'-0.81%"

$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0005981"
$ws.Range("E15").Value = "'-94.18%"

$ws.Range("D16").Value = "'0.005932"
$ws.Range("E16").Value = "'-2.53%"

$ws.Range("D17").Value = "'3.496"
$ws.Range("E17").Value = "'0.04%"

$ws.Range("E18").Value = "'-3.71%"

$ws.Range("E19").Value = "'-0.55%"

$ws.Range("D20").Value = "'0.03359"
$ws.Range("E20").Value = "'1.55%"

$ws.Range("D21").Value = "'0.1307"
$ws.Range("E21").Value = "'1.30%"

$ws.Range("D22").Value = "'3.520"
$ws.Range("E22").Value = "'1.61%"

$ws.Range("D23").Value = "'0.04070"
$ws.Range("E23").Value = "'-1.58%"

$ws.Range("D24").Value = "'0.1378"
$ws.Range("E24").Value = "'-0.06%"

$ws.Range("D25").Value = "'0.001224"
$ws.Range("E25").Value = "'-0.16%"

$ws.Range("E26").Value = "'-16.86%"

$ws.Range("E27").Value = "'-0.75%"

$ws.Range("D40").Value = "'0.03786"
$ws.Range("E40").Value = "'0.79%"

$ws.Range("D41").Value = "'0.1068"
$ws.Range("E41").Value = "'-0.23%"

$ws.Range("D42").Value = "'0.003738"
$ws.Range("E42").Value = "'-35.62%"

$ws.Range("D43").Value = "'0.002401"
$ws.Range("E43").Value = "'4.36%"

$ws.Range("D44").Value = "'0.009468"
$ws.Range("E44").Value = "'0.01%"

$ws.Range("E45").Value = "'0.34%"

$ws.Range("E46").Value = "'0.07%"

$ws.Range("D47").Value = "'0.08911"
$ws.Range("E47").Value = "'56.42%"

$ws.Range("D48").Value = "'0.002264"
$ws.Range("E48").Value = "'0.20%"

$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'0.07%"

$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'0.07%"
